# WebAutomationTestData.xlsx — "Fixed The Excel Format"
#
# The sheet originally had 4 columns (A:D) x 3 rows, with column A holding
# "TestCase" (row1) / "UserPermission" (row2) labels that don't belong with
# the rest of the table (UserName/Password/... in columns B:D).
#
# The fix removes that stray column A (dropping the now-orphaned "TestCase"
# string entirely) and promotes "UserPermission" to its own header row above
# the remaining 3-column table, which shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column A ("TestCase" / "UserPermission" / blank) — B:D shift left to A:C.
$ws.Columns("A").Delete()

# Insert a fresh row 1 and give "UserPermission" its own single-cell header row.
$ws.Rows("1").Insert()
$ws.Range("A1").Value = "UserPermission"

# Match the author's final selection (cell A2, the new top-left data cell).
$ws.Range("A2").Select()
